# Apply hybrid bold + color (2C3E50) highlighting to quantitative metrics
# in specific bullet paragraphs, matching the target diff.
#
# Word COM Font.Color expects a BGR-packed integer (0x00BBGGRR), so RGB
# 2C3E50 (R=2C G=3E B=50) becomes 0x503E2C = 5258796.
$d = $word.ActiveDocument
$metricColor = 5258796

function Apply-MetricBolding {
    param(
        [int]$ParaIndex,
        [string[]]$Metrics
    )

    $p = $d.Paragraphs.Item($ParaIndex)
    $paraRange = $p.Range

    $cursor = $paraRange.Duplicate
    $cursor.Start = $paraRange.Start
    $cursor.End = $paraRange.End

    foreach ($metric in $Metrics) {
        $search = $cursor.Duplicate
        $search.Start = $cursor.Start
        $search.End = $paraRange.End

        $found = $search.Find.Execute($metric, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

        if ($found) {
            $search.Font.Bold = 1
            $search.Font.Color = $metricColor
            $cursor.Start = $search.End
            $cursor.End = $paraRange.End
        }
    }
}

# • Discovered systematic race coding errors ... accuracy from 23% to 64%
Apply-MetricBolding 9 @("23%", "64%")

# • Achieved 87% prediction accuracy ... standard of 71%, reducing polling
#   error margins from ±4.2% to ±2.1%
Apply-MetricBolding 11 @("87%", "71%", "±4.2%", "±2.1%")

# • Wrote RFP and analyzed bids from 1,200 vendors ...
Apply-MetricBolding 31 @("1,200")

# • Created comprehensive meta-analysis framework ... became the $400M
#   Polling Consortium Database at The Analyst Institute, now valued at $1B+
Apply-MetricBolding 46 @("$400M", "$1B")

# • Algorithm reduced mapping costs by 73.5%, saving campaigns and
#   organizations $4.7M
Apply-MetricBolding 63 @("73.5%", "$4.7M")

# • Achieved 87% prediction accuracy for voter turnout vs. industry
#   standard of 71%
Apply-MetricBolding 65 @("87%", "71%")

Write-Output "done"
